$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("2023-12-13 20:27:30", 0.0014),
    @("2023-12-13 20:28:02", 0.0022),
    @("2023-12-13 20:28:27", 0.0014),
    @("2023-12-13 20:28:33", 0.0004)
)

$startRow = 298
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
